$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-removed ECs-as-sender rows (old rows 8,9,10) -- delete bottom-up
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()

# Refresh TPM-derived values for remaining rows 2-7
# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Angpt1"
$ws.Cells.Item(2, 3).Value = "Tek"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 11.05178533333333
$ws.Cells.Item(2, 8).Value = 33.155356
$ws.Cells.Item(2, 9).Value = 0.9017494976312432
$ws.Cells.Item(2, 10).Value = 0.9017494976312432
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 43.28265566666667
$ws.Cells.Item(2, 14).Value = 129.847967
$ws.Cells.Item(2, 15).Value = 0.667219228070094
$ws.Cells.Item(2, 16).Value = 0.667219228070094
$ws.Cells.Item(2, 17).Value = 478.3506190845836
$ws.Cells.Item(2, 18).Value = 4305.155571761252
$ws.Cells.Item(2, 19).Value = 0.6016646037221132
$ws.Cells.Item(2, 20).Value = 0.6016646037221132

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Angpt1"
$ws.Cells.Item(3, 3).Value = "Tek"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 11.05178533333333
$ws.Cells.Item(3, 8).Value = 33.155356
$ws.Cells.Item(3, 9).Value = 0.9017494976312432
$ws.Cells.Item(3, 10).Value = 0.9017494976312432
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 18.63243533333333
$ws.Cells.Item(3, 14).Value = 55.897306
$ws.Cells.Item(3, 15).Value = 0.2872263480299067
$ws.Cells.Item(3, 16).Value = 0.2872263480299067
$ws.Cells.Item(3, 17).Value = 205.9216755412151
$ws.Cells.Item(3, 18).Value = 1853.295079870936
$ws.Cells.Item(3, 19).Value = 0.2590062150424249
$ws.Cells.Item(3, 20).Value = 0.2590062150424249

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Angpt1"
$ws.Cells.Item(4, 3).Value = "Tek"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 11.05178533333333
$ws.Cells.Item(4, 8).Value = 33.155356
$ws.Cells.Item(4, 9).Value = 0.9017494976312432
$ws.Cells.Item(4, 10).Value = 0.9017494976312432
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.955125333333334
$ws.Cells.Item(4, 14).Value = 8.865376000000001
$ws.Cells.Item(4, 15).Value = 0.04555442389999943
$ws.Cells.Item(4, 16).Value = 0.04555442389999944
$ws.Cells.Item(4, 17).Value = 32.65941081709511
$ws.Cells.Item(4, 18).Value = 293.934697353856
$ws.Cells.Item(4, 19).Value = 0.04107867886670519
$ws.Cells.Item(4, 20).Value = 0.04107867886670519

# Row 5
$ws.Cells.Item(5, 1).Value = "MuSCs"
$ws.Cells.Item(5, 2).Value = "Angpt1"
$ws.Cells.Item(5, 3).Value = "Tek"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.204152
$ws.Cells.Item(5, 8).Value = 3.612456
$ws.Cells.Item(5, 9).Value = 0.09825050236875665
$ws.Cells.Item(5, 10).Value = 0.09825050236875667
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 43.28265566666667
$ws.Cells.Item(5, 14).Value = 129.847967
$ws.Cells.Item(5, 15).Value = 0.667219228070094
$ws.Cells.Item(5, 16).Value = 0.667219228070094
$ws.Cells.Item(5, 17).Value = 52.118896386328
$ws.Cells.Item(5, 18).Value = 469.070067476952
$ws.Cells.Item(5, 19).Value = 0.06555462434798076
$ws.Cells.Item(5, 20).Value = 0.06555462434798077

# Row 6
$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 2).Value = "Angpt1"
$ws.Cells.Item(6, 3).Value = "Tek"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.204152
$ws.Cells.Item(6, 8).Value = 3.612456
$ws.Cells.Item(6, 9).Value = 0.09825050236875665
$ws.Cells.Item(6, 10).Value = 0.09825050236875667
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 18.63243533333333
$ws.Cells.Item(6, 14).Value = 55.897306
$ws.Cells.Item(6, 15).Value = 0.2872263480299067
$ws.Cells.Item(6, 16).Value = 0.2872263480299067
$ws.Cells.Item(6, 17).Value = 22.436284271504
$ws.Cells.Item(6, 18).Value = 201.926558443536
$ws.Cells.Item(6, 19).Value = 0.02822013298748167
$ws.Cells.Item(6, 20).Value = 0.02822013298748167

# Row 7
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Angpt1"
$ws.Cells.Item(7, 3).Value = "Tek"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.204152
$ws.Cells.Item(7, 8).Value = 3.612456
$ws.Cells.Item(7, 9).Value = 0.09825050236875665
$ws.Cells.Item(7, 10).Value = 0.09825050236875667
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.955125333333334
$ws.Cells.Item(7, 14).Value = 8.865376000000001
$ws.Cells.Item(7, 15).Value = 0.04555442389999943
$ws.Cells.Item(7, 16).Value = 0.04555442389999944
$ws.Cells.Item(7, 17).Value = 3.558420080384
$ws.Cells.Item(7, 18).Value = 32.025780723456
$ws.Cells.Item(7, 19).Value = 0.004475745033294239
$ws.Cells.Item(7, 20).Value = 0.00447574503329424
